# "run 8 and run 9 E5 titrations and blue tank titrations"
# Append the new CRM-accuracy measurement (run taken 2021-07-06) as row 44,
# following the same layout as the existing rows (Date, CRM value,
# Batch value, % off formula, Batch #, Notes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 44

$ws.Range("A$newRow").Value = 20210706
$ws.Range("B$newRow").Value = 2230.0918634412801
$ws.Range("C$newRow").Value = 2224.4699999999998
$ws.Range("D$newRow").Formula = "=100*(B$newRow-C$newRow)/C$newRow"
$ws.Range("E$newRow").Value = 180
$ws.Range("F$newRow").Value = "CRM OPENED 20210706"

# Matches the author's last selection after entering the new row.
$ws.Range("F$newRow").Select()
